{"js": "// Replace the 100 arithmetic-equation cell values in the single table with\n// their updated answers, preserving all formatting (fonts, size, alignment,\n// cell widths, etc.) \u2014 only the <w:t> text content of each cell changes.\n// Mapping is strictly positional (row-major, top-left to bottom-right) to\n// correctly handle the few equations whose original text repeats.\n\nconst newValues = [\n  [\"42+35=77\", \"73-29=44\", \"11+73=84\", \"50-27=23\", \"75-70=5\"],\n  [\"63-56=7\", \"44+2=46\", \"56-19=37\", \"11+22=33\", \"26-4=22\"],\n  [\"94-72=22\", \"71+17=88\", \"16+78=94\", \"21+36=57\", \"57+16=73\"],\n  [\"61+20=81\", \"94-44=50\", \"45+16=61\", \"97+2=99\", \"20+69=89\"],\n  [\"37+9=46\", \"86-84=2\", \"87-81=6\", \"43-39=4\", \"6+10=16\"],\n  [\"85-46=39\", \"5+75=80\", \"83+14=97\", \"4+67=71\", \"28+11=39\"],\n  [\"20+74=94\", \"20-15=5\", \"52-15=37\", \"23-2=21\", \"67+16=83\"],\n  [\"83-1=82\", \"85-3=82\", \"93-85=8\", \"54-48=6\", \"93-27=66\"],\n  [\"17+1=18\", \"57-50=7\", \"1+31=32\", \"40-31=9\", \"76-63=13\"],\n  [\"63+22=85\", \"81+12=93\", \"55-16=39\", \"23+21=44\", \"95-54=41\"],\n  [\"38+16=54\", \"31-12=19\", \"71-9=62\", \"79-58=21\", \"10+77=87\"],\n  [\"48+0=48\", \"39-0=39\", \"58-4=54\", \"54-45=9\", \"25-4=21\"],\n  [\"90-68=22\", \"65+26=91\", \"97-21=76\", \"81-63=18\", \"21+72=93\"],\n  [\"24+66=90\", \"16+6=22\", \"54+42=96\", \"58-3=55\", \"37-37=0\"],\n  [\"40-19=21\", \"42+22=64\", \"96-24=72\", \"60-0=60\", \"79-56=23\"],\n  [\"89-39=50\", \"14+39=53\", \"92-54=38\", \"64-45=19\", \"53+18=71\"],\n  [\"26+42=68\", \"62-14=48\", \"89-76=13\", \"51+15=66\", \"49+48=97\"],\n  [\"20+19=39\", \"1+49=50\", \"97-78=19\", \"16+66=82\", \"30+22=52\"],\n  [\"18+3=21\", \"59-27=32\", \"25+26=51\", \"35+63=98\", \"77+22=99\"],\n  [\"24+61=85\", \"74-14=60\", \"81-67=14\", \"52-47=5\", \"9+51=60\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nif (table.rowCount !== newValues.length) {\n  throw new Error(\n    \"Unexpected row count: \" + table.rowCount + \" vs expected \" + newValues.length\n  );\n}\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Replace the 100 arithmetic-equation cell values in the single table with\n# their updated answers, preserving all formatting (fonts, size, alignment,\n# cell widths, etc.) - only the text of each cell changes.\n# Mapping is strictly positional (row-major, top-left to bottom-right) to\n# correctly handle the few equations whose original text repeats.\n\n$newValues = @(\n    @(\"42+35=77\", \"73-29=44\", \"11+73=84\", \"50-27=23\", \"75-70=5\"),\n    @(\"63-56=7\", \"44+2=46\", \"56-19=37\", \"11+22=33\", \"26-4=22\"),\n    @(\"94-72=22\", \"71+17=88\", \"16+78=94\", \"21+36=57\", \"57+16=73\"),\n    @(\"61+20=81\", \"94-44=50\", \"45+16=61\", \"97+2=99\", \"20+69=89\"),\n    @(\"37+9=46\", \"86-84=2\", \"87-81=6\", \"43-39=4\", \"6+10=16\"),\n    @(\"85-46=39\", \"5+75=80\", \"83+14=97\", \"4+67=71\", \"28+11=39\"),\n    @(\"20+74=94\", \"20-15=5\", \"52-15=37\", \"23-2=21\", \"67+16=83\"),\n    @(\"83-1=82\", \"85-3=82\", \"93-85=8\", \"54-48=6\", \"93-27=66\"),\n    @(\"17+1=18\", \"57-50=7\", \"1+31=32\", \"40-31=9\", \"76-63=13\"),\n    @(\"63+22=85\", \"81+12=93\", \"55-16=39\", \"23+21=44\", \"95-54=41\"),\n    @(\"38+16=54\", \"31-12=19\", \"71-9=62\", \"79-58=21\", \"10+77=87\"),\n    @(\"48+0=48\", \"39-0=39\", \"58-4=54\", \"54-45=9\", \"25-4=21\"),\n    @(\"90-68=22\", \"65+26=91\", \"97-21=76\", \"81-63=18\", \"21+72=93\"),\n    @(\"24+66=90\", \"16+6=22\", \"54+42=96\", \"58-3=55\", \"37-37=0\"),\n    @(\"40-19=21\", \"42+22=64\", \"96-24=72\", \"60-0=60\", \"79-56=23\"),\n    @(\"89-39=50\", \"14+39=53\", \"92-54=38\", \"64-45=19\", \"53+18=71\"),\n    @(\"26+42=68\", \"62-14=48\", \"89-76=13\", \"51+15=66\", \"49+48=97\"),\n    @(\"20+19=39\", \"1+49=50\", \"97-78=19\", \"16+66=82\", \"30+22=52\"),\n    @(\"18+3=21\", \"59-27=32\", \"25+26=51\", \"35+63=98\", \"77+22=99\"),\n    @(\"24+61=85\", \"74-14=60\", \"81-67=14\", \"52-47=5\", \"9+51=60\")\n)\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\nif ($tbl.Rows.Count -ne $newValues.Count) {\n    throw \"Unexpected row count: $($tbl.Rows.Count) vs expected $($newValues.Count)\"\n}\n\nfor ($r = 1; $r -le $tbl.Rows.Count; $r++) {\n    for ($c = 1; $c -le $tbl.Columns.Count; $c++) {\n        $cell = $tbl.Cell($r, $c)\n        $cell.Range.Text = $newValues[$r - 1][$c - 1]\n    }\n}\n"}
